# Fruta / hortaliza, semanal
#
# A new weekly price observation (Vega Monumental Concepción - Repollo,
# Crespo record, Primera) is inserted as a new data row at row 203,
# pushing the existing rows 203-226 down to 204-227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 203 (shifts 203..226 -> 204..227,
# mirrors Excel's own Rows.Insert which also grows the sheet dimension).
$ws.Rows.Item(203).Insert()

$ws.Cells.Item(203, 1).Value  = 11
$ws.Cells.Item(203, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(203, 3).Value  = "Bíobío"
$ws.Cells.Item(203, 4).Value  = 44505
$ws.Cells.Item(203, 5).Value  = 8
$ws.Cells.Item(203, 6).Value  = 100112006
$ws.Cells.Item(203, 7).Value  = "Repollo"
$ws.Cells.Item(203, 8).Value  = "Crespo record"
$ws.Cells.Item(203, 9).Value  = "Primera"
$ws.Cells.Item(203, 10).Value = 2700
$ws.Cells.Item(203, 11).Value = 650
$ws.Cells.Item(203, 12).Value = 700
$ws.Cells.Item(203, 13).Value = 672
$ws.Cells.Item(203, 14).Value = "$/unidad"
$ws.Cells.Item(203, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(203, 16).Value = 672
$ws.Cells.Item(203, 17).Value = 1
$ws.Cells.Item(203, 18).Value = "Hortaliza"
